$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("queries")

$ws.Range("A6").Value = "verifyModelNodePropertiesEmpty"
$ws.Range("B6").Value = "MATCH (n:node { model: `$modelHandle, version: `$versionString, handle: `$nodeHandle })
`t`t`tMATCH (n)-[:has_property]->(p:property)
`t`t`tRETURN p"

$ws.Range("B6").WrapText = $true
$ws.Rows.Item(6).RowHeight = 51

$ws.Range("A7").Select()
